$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet3")
$ws.Range("A1").Value = "क्रमसं."
$ws.Range("A2").Value = 1

$a1font = $ws.Range("A1").Font
$a1font.Bold = $true; $a1font.Size = 9; $a1font.Color = 0
$ws.Range("A1").WrapText = $true
$ws.Range("A1").HorizontalAlignment = -4108
$ws.Range("A1").VerticalAlignment = -4108

$a2font = $ws.Range("A2").Font
$a2font.Bold = $true; $a2font.Size = 9
$ws.Range("A2").HorizontalAlignment = -4108
$ws.Range("A2").VerticalAlignment = -4108
